$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) data for columns B:AD on rows 172-174.
# Column A ("id") stays fixed per row; only B..AD rotate between the rows.
$row172 = $ws.Range("B172:AD172").Value2
$row173 = $ws.Range("B173:AD173").Value2
$row174 = $ws.Range("B174:AD174").Value2

# Rotate: new row172 <- old row174, new row173 <- old row172, new row174 <- old row173
$ws.Range("B172:AD172").Value2 = $row174
$ws.Range("B173:AD173").Value2 = $row172
$ws.Range("B174:AD174").Value2 = $row173
